$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Sheet "Forecast Comparison": refreshed forecast data (rows 2-17) ---
$ws1.Range("B2").Value = "'2025-02-02"
$ws1.Range("B2").Style = "Normal"
$ws1.Range("D2").Value = 188
$ws1.Range("E2").Value = 145
$ws1.Range("F2").Value = 166
$ws1.Range("G2").Value = 183
$ws1.Range("H2").Value = 208

$ws1.Range("B3").Value = "'2025-02-09"
$ws1.Range("B3").Style = "Normal"
$ws1.Range("D3").Value = 183
$ws1.Range("E3").Value = 141
$ws1.Range("F3").Value = 165
$ws1.Range("G3").Value = 185
$ws1.Range("H3").Value = 216

$ws1.Range("B4").Value = "'2025-02-16"
$ws1.Range("B4").Style = "Normal"
$ws1.Range("D4").Value = 179
$ws1.Range("E4").Value = 138
$ws1.Range("F4").Value = 162
$ws1.Range("G4").Value = 185
$ws1.Range("H4").Value = 220

$ws1.Range("B5").Value = "'2025-02-23"
$ws1.Range("B5").Style = "Normal"
$ws1.Range("D5").Value = 179
$ws1.Range("E5").Value = 138
$ws1.Range("F5").Value = 163
$ws1.Range("G5").Value = 186
$ws1.Range("H5").Value = 221

$ws1.Range("B6").Value = "'2025-03-02"
$ws1.Range("B6").Style = "Normal"
$ws1.Range("D6").Value = 183
$ws1.Range("E6").Value = 141
$ws1.Range("F6").Value = 167
$ws1.Range("G6").Value = 191
$ws1.Range("H6").Value = 230

$ws1.Range("B7").Value = "'2025-03-09"
$ws1.Range("B7").Style = "Normal"
$ws1.Range("D7").Value = 174
$ws1.Range("E7").Value = 134
$ws1.Range("F7").Value = 160
$ws1.Range("G7").Value = 185
$ws1.Range("H7").Value = 224

$ws1.Range("B8").Value = "'2025-03-16"
$ws1.Range("B8").Style = "Normal"
$ws1.Range("D8").Value = 178
$ws1.Range("E8").Value = 137
$ws1.Range("F8").Value = 164
$ws1.Range("G8").Value = 192
$ws1.Range("H8").Value = 236

$ws1.Range("B9").Value = "'2025-03-23"
$ws1.Range("B9").Style = "Normal"
$ws1.Range("D9").Value = 176
$ws1.Range("E9").Value = 135
$ws1.Range("F9").Value = 163
$ws1.Range("G9").Value = 190
$ws1.Range("H9").Value = 234

$ws1.Range("B10").Value = "'2025-03-30"
$ws1.Range("B10").Style = "Normal"
$ws1.Range("D10").Value = 173
$ws1.Range("E10").Value = 133
$ws1.Range("F10").Value = 159
$ws1.Range("G10").Value = 184
$ws1.Range("H10").Value = 224

$ws1.Range("B11").Value = "'2025-04-06"
$ws1.Range("B11").Style = "Normal"
$ws1.Range("D11").Value = 170
$ws1.Range("E11").Value = 131
$ws1.Range("F11").Value = 157
$ws1.Range("G11").Value = 185
$ws1.Range("H11").Value = 227

$ws1.Range("B12").Value = "'2025-04-13"
$ws1.Range("B12").Style = "Normal"
$ws1.Range("D12").Value = 169
$ws1.Range("E12").Value = 130
$ws1.Range("F12").Value = 157
$ws1.Range("G12").Value = 186
$ws1.Range("H12").Value = 231

$ws1.Range("B13").Value = "'2025-04-20"
$ws1.Range("B13").Style = "Normal"
$ws1.Range("D13").Value = 165
$ws1.Range("E13").Value = 127
$ws1.Range("F13").Value = 153
$ws1.Range("G13").Value = 181
$ws1.Range("H13").Value = 223

$ws1.Range("B14").Value = "'2025-04-27"
$ws1.Range("B14").Style = "Normal"
$ws1.Range("D14").Value = 164
$ws1.Range("E14").Value = 126
$ws1.Range("F14").Value = 151
$ws1.Range("G14").Value = 177
$ws1.Range("H14").Value = 218

$ws1.Range("B15").Value = "'2025-05-04"
$ws1.Range("B15").Style = "Normal"
$ws1.Range("D15").Value = 159
$ws1.Range("E15").Value = 122
$ws1.Range("F15").Value = 147
$ws1.Range("G15").Value = 174
$ws1.Range("H15").Value = 216

$ws1.Range("B16").Value = "'2025-05-11"
$ws1.Range("B16").Style = "Normal"
$ws1.Range("D16").Value = 156
$ws1.Range("E16").Value = 120
$ws1.Range("F16").Value = 144
$ws1.Range("G16").Value = 170
$ws1.Range("H16").Value = 209

$ws1.Range("B17").Value = "'2025-05-18"
$ws1.Range("B17").Style = "Normal"
$ws1.Range("D17").Value = 160
$ws1.Range("E17").Value = 123
$ws1.Range("F17").Value = 149
$ws1.Range("G17").Value = 177
$ws1.Range("H17").Value = 221

# --- Sheet "Summary": updated headline metrics ---
$ws2.Range("B2").Value = "'2023-01-01 to 2025-01-26"
$ws2.Range("B2").Style = "Normal"
$ws2.Range("B4").Value = "'405"
$ws2.Range("B4").Style = "Normal"
$ws2.Range("B5").Value = "'131"
$ws2.Range("B5").Style = "Normal"
$ws2.Range("B6").Value = "'114"
$ws2.Range("B6").Style = "Normal"
$ws2.Range("B7").Value = "'86"
$ws2.Range("B7").Style = "Normal"
$ws2.Range("B8").Value = "'13969 units"
$ws2.Range("B8").Style = "Normal"
$ws2.Range("B9").Value = "'2757"
$ws2.Range("B9").Style = "Normal"
$ws2.Range("B10").Value = "'1442"
$ws2.Range("B10").Style = "Normal"
$ws2.Range("B11").Value = "'731"
$ws2.Range("B11").Style = "Normal"
$ws2.Range("B12").Value = "'188"
$ws2.Range("B12").Style = "Normal"
$ws2.Range("B13").Value = "'2025-02-02"
$ws2.Range("B13").Style = "Normal"
$ws2.Range("B14").Value = "'156"
$ws2.Range("B14").Style = "Normal"
